# Applies scheduled market-price / profit recalculation updates to the
# Behemoth_Profits workbook (one leve-profit table per job sheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 295
$ws.Range("I9").Value = 344.2
$ws.Range("J9").Value = 233.5
$ws.Range("K9").Value = 344.2
$ws.Range("L9").Value = 233.5
$ws.Range("M9").Value = -175.2
$ws.Range("N9").Value = -571.5
# Row 34
$ws.Range("H34").Value = 51666
$ws.Range("I34").Value = 51666
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 51666
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -51463
# Row 36
$ws.Range("H36").Value = 51666
$ws.Range("I36").Value = 51666
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 51666
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -50951
# Row 80
$ws.Range("H80").Value = 1701.5264
$ws.Range("I80").Value = 2323
$ws.Range("J80").Value = 636.1429000000001
$ws.Range("K80").Value = 6969
$ws.Range("L80").Value = 1908.4287
$ws.Range("M80").Value = -5971
$ws.Range("N80").Value = -3904.4287
# Row 83
$ws.Range("H83").Value = 1701.5264
$ws.Range("I83").Value = 2323
$ws.Range("J83").Value = 636.1429000000001
$ws.Range("K83").Value = 20907
$ws.Range("L83").Value = 5725.2861
$ws.Range("M83").Value = -15915
$ws.Range("N83").Value = -15709.2861
# Row 86
$ws.Range("H86").Value = 4615.231
$ws.Range("J86").Value = 5139.9
$ws.Range("L86").Value = 5139.9
$ws.Range("N86").Value = -7385.9
# Row 89
$ws.Range("H89").Value = 4615.231
$ws.Range("J89").Value = 5139.9
$ws.Range("L89").Value = 25699.5
$ws.Range("N89").Value = -36931.5
# Row 92
$ws.Range("H92").Value = 463.1
$ws.Range("I92").Value = 484.82352
$ws.Range("J92").Value = 340
$ws.Range("K92").Value = 484.82352
$ws.Range("L92").Value = 340
$ws.Range("M92").Value = 763.1764800000001
$ws.Range("N92").Value = -2836
# Row 106
$ws.Range("H106").Value = 9103.684999999999
$ws.Range("I106").Value = 2139.2856
$ws.Range("K106").Value = 2139.2856
$ws.Range("M106").Value = -1508.2856
# Row 124
$ws.Range("H124").Value = 74979.5
$ws.Range("J124").Value = 74979.5
$ws.Range("L124").Value = 74979.5
$ws.Range("N124").Value = -84799.5
# Row 127
$ws.Range("H127").Value = 8819.299999999999
$ws.Range("I127").Value = 831.3333
$ws.Range("J127").Value = 12242.714
$ws.Range("K127").Value = 2493.9999
$ws.Range("L127").Value = 36728.142
$ws.Range("M127").Value = 2466.0001
$ws.Range("N127").Value = -46648.142
# Row 129
$ws.Range("H129").Value = 3143.7778
$ws.Range("I129").Value = 1758.8
$ws.Range("K129").Value = 5276.4
$ws.Range("M129").Value = -276.3999999999996
# Row 131
$ws.Range("H131").Value = 4780.2383
$ws.Range("I131").Value = 3313.5715
$ws.Range("K131").Value = 9940.7145
$ws.Range("M131").Value = -4900.7145
# Row 137
$ws.Range("H137").Value = 1011416.4
$ws.Range("I137").Value = 2005099
$ws.Range("J137").Value = 17733.8
$ws.Range("K137").Value = 6015297
$ws.Range("L137").Value = 53201.39999999999
$ws.Range("M137").Value = -6012747
$ws.Range("N137").Value = -58301.39999999999
# Row 141
$ws.Range("H141").Value = 7934.5557
$ws.Range("I141").Value = 9682.799999999999
$ws.Range("K141").Value = 29048.4
$ws.Range("M141").Value = -23868.4

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5839586.5
$ws.Range("I32").Value = 8082848.5
$ws.Range("J32").Value = 44492.375
$ws.Range("K32").Value = 8082848.5
$ws.Range("L32").Value = 44492.375
$ws.Range("M32").Value = -8082561.5
$ws.Range("N32").Value = -45066.375
# Row 61
$ws.Range("H61").Value = 30004494
$ws.Range("I61").Value = 38465828
$ws.Range("K61").Value = 38465828
$ws.Range("M61").Value = -38465616
# Row 74
$ws.Range("H74").Value = 13170454
$ws.Range("I74").Value = 25005932
$ws.Range("J74").Value = 19923.223
$ws.Range("K74").Value = 25005932
$ws.Range("L74").Value = 19923.223
$ws.Range("M74").Value = -25005058
$ws.Range("N74").Value = -21671.223
# Row 77
$ws.Range("H77").Value = 13170454
$ws.Range("I77").Value = 25005932
$ws.Range("J77").Value = 19923.223
$ws.Range("K77").Value = 125029660
$ws.Range("L77").Value = 99616.11500000001
$ws.Range("M77").Value = -125025292
$ws.Range("N77").Value = -108352.115
# Row 80
$ws.Range("H80").Value = 22000
$ws.Range("I80").Value = 22000
$ws.Range("K80").Value = 22000
$ws.Range("M80").Value = -21002
# Row 83
$ws.Range("H83").Value = 22000
$ws.Range("I83").Value = 22000
$ws.Range("K83").Value = 66000
$ws.Range("M83").Value = -61008
# Row 88
$ws.Range("H88").Value = 3159.4666
$ws.Range("I88").Value = 2760.625
$ws.Range("J88").Value = 3615.2856
$ws.Range("K88").Value = 2760.625
$ws.Range("L88").Value = 3615.2856
$ws.Range("M88").Value = -2354.625
$ws.Range("N88").Value = -4427.2856
# Row 91
$ws.Range("H91").Value = 3159.4666
$ws.Range("I91").Value = 2760.625
$ws.Range("J91").Value = 3615.2856
$ws.Range("K91").Value = 2760.625
$ws.Range("L91").Value = 3615.2856
$ws.Range("M91").Value = -1356.625
$ws.Range("N91").Value = -6423.2856
# Row 110
$ws.Range("H110").Value = 1133.2142
$ws.Range("I110").Value = 1074.2307
$ws.Range("K110").Value = 1074.2307
$ws.Range("M110").Value = 970.7692999999999
# Row 122
$ws.Range("H122").Value = 3550
$ws.Range("I122").Value = 2680
$ws.Range("K122").Value = 8040
$ws.Range("M122").Value = -5590
# Row 136
$ws.Range("H136").Value = 30004494
$ws.Range("I136").Value = 38465828
$ws.Range("K136").Value = 115397484
$ws.Range("M136").Value = -115394934

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 35848.8
$ws.Range("I82").Value = 21418.572
$ws.Range("J82").Value = 69519.336
$ws.Range("K82").Value = 21418.572
$ws.Range("L82").Value = 69519.336
$ws.Range("M82").Value = -21035.572
$ws.Range("N82").Value = -70285.336
# Row 85
$ws.Range("H85").Value = 35848.8
$ws.Range("I85").Value = 21418.572
$ws.Range("J85").Value = 69519.336
$ws.Range("K85").Value = 21418.572
$ws.Range("L85").Value = 69519.336
$ws.Range("M85").Value = -20092.572
$ws.Range("N85").Value = -72171.336
# Row 134
$ws.Range("H134").Value = 3404390
$ws.Range("I134").Value = 1989.375
$ws.Range("K134").Value = 5968.125
$ws.Range("M134").Value = -3433.125

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 772.7
$ws.Range("I16").Value = 720
$ws.Range("J16").Value = 851.75
$ws.Range("K16").Value = 720
$ws.Range("L16").Value = 851.75
$ws.Range("M16").Value = -433
$ws.Range("N16").Value = -1425.75
# Row 31
$ws.Range("H31").Value = 2606579.8
$ws.Range("I31").Value = 26861.666
$ws.Range("K31").Value = 26861.666
$ws.Range("M31").Value = -26566.666
# Row 34
$ws.Range("H34").Value = 2606579.8
$ws.Range("I34").Value = 26861.666
$ws.Range("K34").Value = 26861.666
$ws.Range("M34").Value = -26659.666
# Row 55
$ws.Range("H55").Value = 21099.6
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 21099.6
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 21099.6
$ws.Range("N55").Value = -21729.6
$ws.Range("M55").ClearContents()
# Row 86
$ws.Range("H86").Value = 50247.273
$ws.Range("I86").Value = 4537.0835
$ws.Range("J86").Value = 105099.5
$ws.Range("K86").Value = 4537.0835
$ws.Range("L86").Value = 105099.5
$ws.Range("M86").Value = -3414.0835
$ws.Range("N86").Value = -107345.5
# Row 89
$ws.Range("H89").Value = 50247.273
$ws.Range("I89").Value = 4537.0835
$ws.Range("J89").Value = 105099.5
$ws.Range("K89").Value = 22685.4175
$ws.Range("L89").Value = 525497.5
$ws.Range("M89").Value = -17069.4175
$ws.Range("N89").Value = -536729.5
# Row 94
$ws.Range("H94").Value = 3741.3076
$ws.Range("I94").Value = 3689.9167
$ws.Range("J94").Value = 3785.3572
$ws.Range("K94").Value = 3689.9167
$ws.Range("L94").Value = 3785.3572
$ws.Range("M94").Value = -3238.9167
$ws.Range("N94").Value = -4687.3572
# Row 99
$ws.Range("H99").Value = 3134.7144
$ws.Range("I99").Value = 2849.75
$ws.Range("J99").Value = 3514.6667
$ws.Range("K99").Value = 2849.75
$ws.Range("L99").Value = 3514.6667
$ws.Range("M99").Value = -1351.75
$ws.Range("N99").Value = -6510.6667
# Row 113
$ws.Range("H113").Value = 772.7
$ws.Range("I113").Value = 720
$ws.Range("J113").Value = 851.75
$ws.Range("K113").Value = 720
$ws.Range("L113").Value = 851.75
$ws.Range("M113").Value = 1450
$ws.Range("N113").Value = -5191.75
# Row 126
$ws.Range("H126").Value = 3134.7144
$ws.Range("I126").Value = 2849.75
$ws.Range("J126").Value = 3514.6667
$ws.Range("K126").Value = 8549.25
$ws.Range("L126").Value = 10544.0001
$ws.Range("M126").Value = -6079.25
$ws.Range("N126").Value = -15484.0001
# Row 134
$ws.Range("H134").Value = 3769.1765
$ws.Range("I134").Value = 2346.9167
$ws.Range("J134").Value = 7182.6
$ws.Range("K134").Value = 7040.750100000001
$ws.Range("L134").Value = 21547.8
$ws.Range("M134").Value = -4505.750100000001
$ws.Range("N134").Value = -26617.8
# Row 141
$ws.Range("H141").Value = 190497.34
$ws.Range("J141").Value = 203515.11
$ws.Range("L141").Value = 203515.11
$ws.Range("N141").Value = -213875.11

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 163.82353
$ws.Range("I2").Value = 64.27273
$ws.Range("J2").Value = 346.33334
$ws.Range("K2").Value = 385.63638
$ws.Range("L2").Value = 2078.00004
$ws.Range("M2").Value = -272.63638
$ws.Range("N2").Value = -2304.00004
# Row 12
$ws.Range("H12").Value = 124.545456
$ws.Range("I12").Value = 144.75
$ws.Range("J12").Value = 113
$ws.Range("K12").Value = 434.25
$ws.Range("L12").Value = 339
$ws.Range("M12").Value = -261.25
$ws.Range("N12").Value = -685
# Row 113
$ws.Range("H113").Value = 1358.2667
$ws.Range("J113").Value = 1841.3334
$ws.Range("L113").Value = 5524.0002
$ws.Range("N113").Value = -9864.0002
# Row 114
$ws.Range("H114").Value = 4545.3335
$ws.Range("I114").Value = 4333.3335
$ws.Range("K114").Value = 13000.0005
$ws.Range("M114").Value = -9746.000499999998
# Row 117
$ws.Range("H117").Value = 1582.5
$ws.Range("I117").Value = 993
$ws.Range("J117").Value = 1936.2
$ws.Range("K117").Value = 2979
$ws.Range("L117").Value = 5808.6
$ws.Range("M117").Value = 463
$ws.Range("N117").Value = -12692.6
# Row 141
$ws.Range("H141").Value = 509163.34
$ws.Range("I141").Value = 607996
$ws.Range("K141").Value = 1823988
$ws.Range("M141").Value = -1818808

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6324.6
$ws.Range("I70").Value = 6031.6
$ws.Range("J70").Value = 6617.6
$ws.Range("K70").Value = 6031.6
$ws.Range("L70").Value = 6617.6
$ws.Range("M70").Value = -5761.6
$ws.Range("N70").Value = -7157.6
# Row 73
$ws.Range("H73").Value = 6324.6
$ws.Range("I73").Value = 6031.6
$ws.Range("J73").Value = 6617.6
$ws.Range("K73").Value = 6031.6
$ws.Range("L73").Value = 6617.6
$ws.Range("M73").Value = -5095.6
$ws.Range("N73").Value = -8489.6
# Row 134
$ws.Range("H134").Value = 17000
$ws.Range("J134").Value = 17000
$ws.Range("L134").Value = 51000
$ws.Range("N134").Value = -56070

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5977.8
$ws.Range("I40").Value = 6710.125
$ws.Range("J40").Value = 5140.857
$ws.Range("K40").Value = 6710.125
$ws.Range("L40").Value = 5140.857
$ws.Range("M40").Value = -6574.125
$ws.Range("N40").Value = -5412.857
# Row 46
$ws.Range("H46").Value = 3820.3157
$ws.Range("I46").Value = 2885.4443
$ws.Range("K46").Value = 2885.4443
$ws.Range("M46").Value = -2697.4443
# Row 50
$ws.Range("H50").Value = 36246.5
$ws.Range("I50").Value = 29998
$ws.Range("K50").Value = 29998
$ws.Range("M50").Value = -29361
# Row 61
$ws.Range("H61").Value = 2156.9473
$ws.Range("I61").Value = 2171.8572
$ws.Range("K61").Value = 2171.8572
$ws.Range("M61").Value = -1969.8572
# Row 93
$ws.Range("H93").Value = 62502190
$ws.Range("I93").Value = 90911110
$ws.Range("J93").Value = 2557.8
$ws.Range("K93").Value = 90911110
$ws.Range("L93").Value = 2557.8
$ws.Range("M93").Value = -90909862
$ws.Range("N93").Value = -5053.8
# Row 113
$ws.Range("H113").Value = 2156.9473
$ws.Range("I113").Value = 2171.8572
$ws.Range("K113").Value = 2171.8572
$ws.Range("M113").Value = -1.857199999999921
# Row 122
$ws.Range("H122").Value = 5517.8096
$ws.Range("I122").Value = 5038.4
$ws.Range("J122").Value = 6716.3335
$ws.Range("K122").Value = 15115.2
$ws.Range("L122").Value = 20149.0005
$ws.Range("M122").Value = -12665.2
$ws.Range("N122").Value = -25049.0005
# Row 132
$ws.Range("H132").Value = 4876.7144
$ws.Range("I132").Value = 5289.5
$ws.Range("J132").Value = 3844.75
$ws.Range("K132").Value = 15868.5
$ws.Range("L132").Value = 11534.25
$ws.Range("M132").Value = -13338.5
$ws.Range("N132").Value = -16594.25

$ws = $wb.Worksheets.Item("WVR")
# Row 47
$ws.Range("H47").Value = 37495
$ws.Range("J47").Value = 37495
$ws.Range("L47").Value = 37495
$ws.Range("N47").Value = -38639
# Row 126
$ws.Range("H126").Value = 1938.9333
$ws.Range("I126").Value = 2010.2222
$ws.Range("J126").Value = 1832
$ws.Range("K126").Value = 6030.6666
$ws.Range("L126").Value = 5496
$ws.Range("M126").Value = -3560.6666
$ws.Range("N126").Value = -10436
# Row 132
$ws.Range("H132").Value = 3793695
$ws.Range("I132").Value = 5067.647
$ws.Range("J132").Value = 16675028
$ws.Range("K132").Value = 15202.941
$ws.Range("L132").Value = 50025084
$ws.Range("M132").Value = -12672.941
$ws.Range("N132").Value = -50030144
# Row 136
$ws.Range("H136").Value = 4427
$ws.Range("I136").Value = 3639.8215
$ws.Range("K136").Value = 10919.4645
$ws.Range("M136").Value = -8369.4645
# Row 139
$ws.Range("H139").Value = 84000
